$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "StudentSheet"

# Update header row (row 1) and the first data row, in the order the
# original author entered them (new shared strings are appended in entry
# order): name, then 3A, then the remaining header labels.
$ws.Range("A1").Value = "name"

$ws.Range("A2").Value = "Cruz, Aimee Lou D."
$ws.Range("B2").Value = "BSIT"
$ws.Range("C2").Value = "3A"
$ws.Range("D2").Value = "2017-2018"

$ws.Range("B1").Value = "class_course"
$ws.Range("C1").Value = "class_section"
$ws.Range("D1").Value = "class_school_year"

# Row 3 / Row 4: only keep the name, clear the rest of the row
$ws.Range("A3").Value = "De Vera, Jazelene Mae M."
$ws.Range("B3:D3").ClearContents()

$ws.Range("A4").Value = "George, Donald Patrick C."
$ws.Range("B4:D4").ClearContents()

# Auto-fit columns to content (bestFit widths)
$ws.Columns("A:D").AutoFit()

# Update selection to column B (whole column)
$ws.Range("B1:B1048576").Select()
